# Update the phone number of the contact person on sheet "6.2.1.1".
# Row 9 holds "Телефон контактного лица" (Phone number of contact person)
# in column A, with the actual phone number in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("6.2.1.1")
$ws.Range("B9").Value = "(312) 32 46 55"
